$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- G9: rename the task description (old string becomes unused and is
#     garbage-collected from sharedStrings; new strings get appended in the
#     order below so the resulting <sst> index order matches the target).
$ws.Range("G14").Value = "Logic Unit Timing Simulation and Timing Waves Screenshots"
$ws.Range("G9").Value  = "Write LogicUnit.vhd + Set up Testing Environment"
$ws.Range("G15").Value = "Fixed Logic of AltB and AltBu"

# --- Row 14: new activity-log entry (8414, 2020-03-30, 21:00 - 21:15)
$ws.Range("B14").Value = 8414
$ws.Range("C14").Value = 43920
$ws.Range("D14").Value = 0.875
$ws.Range("E14").Value = 0.88541666666666663
# D14 picks up the bordered time style (matches E-column's look), like the
# target XML's s="22" for D14.
$ws.Range("E9").Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4122) | Out-Null

# --- Row 15: new activity-log entry (8414, 2020-03-30, 21:30 - 22:30)
#     D15 already carries the right (unbordered) time style from the
#     template row, so only the values need to be written.
$ws.Range("B15").Value = 8414
$ws.Range("C15").Value = 43920
$ws.Range("D15").Value = 0.89583333333333337
$ws.Range("E15").Value = 0.9375

# --- Row 16: new activity-log entry (8414, 2020-03-30, started 22:30)
$ws.Range("B16").Value = 8414
$ws.Range("C16").Value = 43920
$ws.Range("D16").Value = 0.9375
$ws.Range("E9").Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Sheet view: move the viewport up and move the live selection to G10.
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("G10").Select()
